# 6.4.2.1 Total freshwater withdrawal — add the 2023 data column (M) and
# the trilingual "Items" header row labels (A4:C4), matching the upstream
# gh-pages data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New trilingual "Items" column header (row 4, A:C) -------------------
# Setting these in Kyrgyz / Russian / English order appends three new
# shared strings at the end of the table (indices 49, 50, 51), exactly как
# upstream did.
$ws.Range("A4").Value = "Көрсөткүчтөрдүн аталыштары"
$ws.Range("B4").Value = "Наименование показателей"
$ws.Range("C4").Value = "Items"

# --- New year column M: 2023 data -----------------------------------------
$ws.Range("M4").Value = 2023
$ws.Range("M5").Value = 8872.5
$ws.Range("M7").Value = 8601.5
$ws.Range("M8").Value = 271
$ws.Range("M10").Value = 723.4
$ws.Range("M11").Value = 1205.5999999999999
$ws.Range("M12").Value = 779.6
$ws.Range("M13").Value = 829.3
$ws.Range("M14").Value = 1314.9
$ws.Range("M15").Value = 1034.5999999999999
$ws.Range("M16").Value = 2762.1
$ws.Range("M17").Value = 166
$ws.Range("M18").Value = 57

# Rows 6 and 9 are section-header rows whose data cells (D:L) are
# intentionally blank but still carry the row's number formatting; mirror
# that onto the new M column.
$ws.Range("M6").NumberFormat = $ws.Range("L6").NumberFormat
$ws.Range("M9").NumberFormat = $ws.Range("L9").NumberFormat

# --- Carry the column L formatting across into the new column M ----------
# (font, borders, alignment, number format) so the new column matches the
# rest of the 2014-2022 series visually.
$ws.Range("L4:L18").Copy()
$ws.Range("M4:M18").PasteSpecial(-4122)
$ws.Range("M4").Select()
$excel.CutCopyMode = $false

# Re-apply the values (PasteSpecial(xlPasteFormats) only touches
# formatting, but make sure nothing was clobbered).
$ws.Range("M4").Value = 2023
$ws.Range("M5").Value = 8872.5
$ws.Range("M7").Value = 8601.5
$ws.Range("M8").Value = 271
$ws.Range("M10").Value = 723.4
$ws.Range("M11").Value = 1205.5999999999999
$ws.Range("M12").Value = 779.6
$ws.Range("M13").Value = 829.3
$ws.Range("M14").Value = 1314.9
$ws.Range("M15").Value = 1034.5999999999999
$ws.Range("M16").Value = 2762.1
$ws.Range("M17").Value = 166
$ws.Range("M18").Value = 57

# --- Row heights: the whole data block (rows 4-18) now uses a uniform
# 14.25pt custom height (was a mix of default/15/12.75). -------------------
for ($r = 4; $r -le 18; $r++) {
    $ws.Rows.Item($r).RowHeight = 14.25
}

# --- Reset the lingering selection back to the sheet's top-left cell -----
$ws.Range("A1").Select()
